$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:K7 with new values
$values = @(
    @(0.241, 0.19, 0.278, 0.267, 0.121, 0.145, 0.146, 0.157, 0.172, 0.334),
    @(0.292, 0.265, 0.347, 0.314, 0.132, 0.222, 0.227, 0.216, 0.233, 0.531),
    @(0.438, 0.413, 0.474, 0.442, 0.167, 0.365, 0.376, 0.412, 0.404, 0.729),
    @(0.574, 0.54,  0.606, 0.608, 0.233, 0.464, 0.469, 0.546, 0.5,   0.859),
    @(0.67,  0.641, 0.704, 0.657, 0.281, 0.542, 0.5570000000000001, 0.623, 0.5659999999999999, 0.929),
    @(0.739, 0.722, 0.794, 0.769, 0.324, 0.638, 0.636, 0.708, 0.666, 0.971)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = 2 + $j  # column B = 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}

# Delete row 8 (the "50" row) entirely
$ws.Rows.Item(8).Delete()
